$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's results as a new row (row 75) at the bottom of the table.
$newRow = 75

$ws.Cells.Item($newRow, 1).Value = 46024
$ws.Cells.Item($newRow, 2).Value = 166
$ws.Cells.Item($newRow, 3).Value = 177
$ws.Cells.Item($newRow, 4).Value = 168

# Match the date formatting used by the rest of column A.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat
